# The original sheet has a redundant first column (A) that duplicates the
# GENE column (originally F). The edit removes that column entirely, which
# shifts all remaining columns (B:F -> A:E) one position to the left,
# carrying their values, types and formatting (incl. the header row's bold
# style) along with them - exactly how Excel's "Delete Entire Column"
# behaves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").EntireColumn.Delete()
